$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data to row 25: Hours = 2, Completed = new paragraph text
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "Update cover image. Update paragraphes."

# Update the sum formula to include the new row
$ws.Range("B42").Formula = "=SUM(B2:B25)"

# Update the view: scroll so row 14 is the top-left visible row, and move
# the active selection to C34
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C34").Select()
